$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2-7
# from serial date 45174 (2023-09-05) to 45175 (2023-09-06)
foreach ($row in 2..7) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
